$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "state_abb_appeldorn_fortunato"
$ws.Range("B80").Value = "State abbreviation in Appeldorn and Fortunato's data"
$ws.Range("B81").Value = "State name in Appeldorn and Fortunato's data"
$ws.Range("A81").Value = "state_name_appeldorn_fortunato"

$ws.Range("A82").Select() | Out-Null
